$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update derived-column values for rows 2-28 (height 169cm -> 169.7cm recalculation) ---
$ws.Range("H2").Value = 0.61
$ws.Range("H3").Value = 0.61
$ws.Range("J4").Value = 0.02
$ws.Range("K4").Value = 28.7
$ws.Range("H5").Value = 0.57999999999999996
$ws.Range("I5").Value = "IMPROVED"
$ws.Range("J5").Value = 0.01
$ws.Range("K5").Value = 28.8
$ws.Range("I6").Value = "WORSENED"
$ws.Range("J6").Value = 0.01
$ws.Range("K6").Value = 28.5
$ws.Range("H7").Value = 0.59
$ws.Range("I7").Value = "SAME"
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 28.6
$ws.Range("J8").Value = 0.01
$ws.Range("K8").Value = 28.4
$ws.Range("H9").Value = 0.57999999999999996
$ws.Range("I9").Value = "SAME"
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 28.5
$ws.Range("I10").Value = "SAME"
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 28.6
$ws.Range("K11").Value = 28.6
$ws.Range("K12").Value = 28.6
$ws.Range("K13").Value = 28.7
$ws.Range("K14").Value = 29
$ws.Range("K15").Value = 28.9
$ws.Range("H16").Value = 0.57999999999999996
$ws.Range("J16").Value = 0.02
$ws.Range("K16").Value = 28.5
$ws.Range("I17").Value = "SAME"
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 28.3
$ws.Range("K18").Value = 28.4
$ws.Range("K19").Value = 28.8
$ws.Range("K20").Value = 28.7
$ws.Range("K21").Value = 28.5
$ws.Range("H22").Value = 0.57999999999999996
$ws.Range("I22").Value = "IMPROVED"
$ws.Range("J22").Value = 0.01
$ws.Range("K22").Value = 28.6
$ws.Range("I23").Value = "WORSENED"
$ws.Range("J23").Value = 0.01
$ws.Range("K23").Value = 28.7
$ws.Range("K24").Value = 28.9
$ws.Range("K25").Value = 29.1
$ws.Range("H26").Value = 0.56999999999999995
$ws.Range("I26").Value = "IMPROVED"
$ws.Range("J26").Value = 0.01
$ws.Range("K26").Value = 28.7
$ws.Range("I27").Value = "WORSENED"
$ws.Range("J27").Value = 0.01
$ws.Range("K27").Value = 28.6
$ws.Range("H28").Value = 0.56999999999999995
$ws.Range("I28").Value = "IMPROVED"
$ws.Range("J28").Value = 0.01
$ws.Range("K28").Value = 28.7

# --- Append two new rows of tracked data (29, 30) ---
# Copy formatting from the last existing data row so number formats / fonts / alignment match.
$ws.Range("A28:L28").Copy() | Out-Null
$ws.Range("A29:L29").PasteSpecial(-4122) | Out-Null
$ws.Range("A28:L28").Copy() | Out-Null
$ws.Range("A30:L30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A29").Value = 44014
$ws.Range("B29").Value = 99.5
$ws.Range("C29").Value = 105.5
$ws.Range("D29").Value = 0.94
$ws.Range("E29").Value = "WORSENED"
$ws.Range("F29").Value = 0.01
$ws.Range("G29").Value = 82.4
$ws.Range("H29").Value = 0.59
$ws.Range("I29").Value = "WORSENED"
$ws.Range("J29").Value = 0.02
$ws.Range("K29").Value = 28.6
$ws.Range("L29").Value = "OVERWEIGHT"

$ws.Range("A30").Value = 44015
$ws.Range("B30").Value = 98.5
$ws.Range("C30").Value = 104.5
$ws.Range("D30").Value = 0.94
$ws.Range("E30").Value = "SAME"
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 81.5
$ws.Range("H30").Value = 0.57999999999999996
$ws.Range("I30").Value = "IMPROVED"
$ws.Range("J30").Value = 0.01
$ws.Range("K30").Value = 28.3
$ws.Range("L30").Value = "OVERWEIGHT"

# --- Remove the no-longer-needed SQL sheet ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("SQL").Delete()
$excel.DisplayAlerts = $true

# --- Rename the data sheet to the plain "운동기록" ---
$ws.Name = "운동기록"

# --- Selection moved to C2 (matches the author's last cursor position) ---
$ws.Range("C2").Select()
